$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep its content as literal text even when the
    # string looks like a plain number (e.g. "572.78"), matching the
    # inlineStr cells already present in the workbook.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.777.66"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.173.43"
$ws.Range("E3").Value = "  -4.63%  "

# Row 5 - BNB
Set-TextValue "D5" "572.78"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6 - Solana
Set-TextValue "D6" "172.26"
$ws.Range("E6").Value = "  -3.43%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.77%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.172.36"
$ws.Range("E9").Value = "  -4.63%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.83%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -3.61%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -3.48%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.722.36"
$ws.Range("E13").Value = "  -4.72%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.78%  "

# Row 15 - Avalanche
Set-TextValue "D15" "27.54"
$ws.Range("E15").Value = "  -3.92%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.747.71"
$ws.Range("E16").Value = "  +0.46%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -2.21%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.175.59"
$ws.Range("E18").Value = "  -5.92%  "

# Row 19 - Polkadot
Set-TextValue "D19" "5.75"
$ws.Range("E19").Value = "  +0.39%  "

# Row 20 - Chainlink
Set-TextValue "D20" "12.96"
$ws.Range("E20").Value = "  -3.24%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "361.74"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.29"
$ws.Range("E22").Value = "  -1.68%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.14%  "

# Row 24 - Litecoin
Set-TextValue "D24" "69.07"
$ws.Range("E24").Value = "  -3.13%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  -4.42%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.308.50"
$ws.Range("E26").Value = "  -4.93%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -6.05%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "9.86"
$ws.Range("E28").Value = "  +3.30%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -0.71%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.05%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.54%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  -0.11%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "5.43"
$ws.Range("E33").Value = "  -3.00%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "22.14"
$ws.Range("E34").Value = "  -3.19%  "

# Row 35 - was Fetch.AI, becomes Aptos
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "6.65"
$ws.Range("E35").Value = "  -2.41%  "

# Row 36 - was Aptos, becomes Fetch.AI
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "1.20"
$ws.Range("E36").Value = "  -0.59%  "

# Row 37 - Monero
Set-TextValue "D37" "160.40"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.44%  "

# Row 39 - Mantle
Set-TextValue "D39" "0.838"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.80"
$ws.Range("E40").Value = "  +3.23%  "

# Row 41 - EnergySwap
Set-TextValue "D41" "26.44"
$ws.Range("E41").Value = "  -2.88%  "

# Row 42 - dogwifhat
Set-TextValue "D42" "2.49"
$ws.Range("E42").Value = "  -3.08%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.645.86"
$ws.Range("E43").Value = "  -2.50%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  -0.94%  "

# Row 45 - Filecoin
Set-TextValue "D45" "4.21"
$ws.Range("E45").Value = "  -1.38%  "

# Row 46 - OKB
Set-TextValue "D46" "39.76"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47 - Bittensor
Set-TextValue "D47" "331.34"
$ws.Range("E47").Value = "  -1.69%  "

# Row 48 - Hedera
Set-TextValue "D48" "0.0656"
$ws.Range("E48").Value = "  -1.47%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "24.04"
$ws.Range("E49").Value = "  +0.35%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -1.10%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -1.10%  "
